$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "92-53=39"
$t.Cell(1,2).Range.Text = "35-23=12"
$t.Cell(1,3).Range.Text = "0+53=53"
$t.Cell(1,4).Range.Text = "75-60=15"
$t.Cell(1,5).Range.Text = "75-68=7"
$t.Cell(2,1).Range.Text = "21+53=74"
$t.Cell(2,2).Range.Text = "91-45=46"
$t.Cell(2,3).Range.Text = "11+70=81"
$t.Cell(2,4).Range.Text = "82-69=13"
$t.Cell(2,5).Range.Text = "67-0=67"
$t.Cell(3,1).Range.Text = "54+15=69"
$t.Cell(3,2).Range.Text = "46-5=41"
$t.Cell(3,3).Range.Text = "51-49=2"
$t.Cell(3,4).Range.Text = "60-52=8"
$t.Cell(3,5).Range.Text = "68+18=86"
$t.Cell(4,1).Range.Text = "90-40=50"
$t.Cell(4,2).Range.Text = "69-65=4"
$t.Cell(4,3).Range.Text = "25-0=25"
$t.Cell(4,4).Range.Text = "93-57=36"
$t.Cell(4,5).Range.Text = "62-21=41"
$t.Cell(5,1).Range.Text = "54-11=43"
$t.Cell(5,2).Range.Text = "57-24=33"
$t.Cell(5,3).Range.Text = "53+24=77"
$t.Cell(5,4).Range.Text = "31+38=69"
$t.Cell(5,5).Range.Text = "5+62=67"
$t.Cell(6,1).Range.Text = "68-36=32"
$t.Cell(6,2).Range.Text = "75-7=68"
$t.Cell(6,3).Range.Text = "46+5=51"
$t.Cell(6,4).Range.Text = "8+16=24"
$t.Cell(6,5).Range.Text = "69-55=14"
$t.Cell(7,1).Range.Text = "49-15=34"
$t.Cell(7,2).Range.Text = "82-79=3"
$t.Cell(7,3).Range.Text = "54+33=87"
$t.Cell(7,4).Range.Text = "8+8=16"
$t.Cell(7,5).Range.Text = "11+29=40"
$t.Cell(8,1).Range.Text = "97-48=49"
$t.Cell(8,2).Range.Text = "37-15=22"
$t.Cell(8,3).Range.Text = "70-39=31"
$t.Cell(8,4).Range.Text = "9+29=38"
$t.Cell(8,5).Range.Text = "30+23=53"
$t.Cell(9,1).Range.Text = "42+47=89"
$t.Cell(9,2).Range.Text = "94-88=6"
$t.Cell(9,3).Range.Text = "65-65=0"
$t.Cell(9,4).Range.Text = "16+63=79"
$t.Cell(9,5).Range.Text = "1+76=77"
$t.Cell(10,1).Range.Text = "59+0=59"
$t.Cell(10,2).Range.Text = "32+57=89"
$t.Cell(10,3).Range.Text = "58-44=14"
$t.Cell(10,4).Range.Text = "49+19=68"
$t.Cell(10,5).Range.Text = "53+43=96"
$t.Cell(11,1).Range.Text = "1+9=10"
$t.Cell(11,2).Range.Text = "14+14=28"
$t.Cell(11,3).Range.Text = "24+62=86"
$t.Cell(11,4).Range.Text = "78-6=72"
$t.Cell(11,5).Range.Text = "87-16=71"
$t.Cell(12,1).Range.Text = "51-1=50"
$t.Cell(12,2).Range.Text = "82-66=16"
$t.Cell(12,3).Range.Text = "84-45=39"
$t.Cell(12,4).Range.Text = "62-44=18"
$t.Cell(12,5).Range.Text = "30+57=87"
$t.Cell(13,1).Range.Text = "81-48=33"
$t.Cell(13,2).Range.Text = "72-31=41"
$t.Cell(13,3).Range.Text = "78-11=67"
$t.Cell(13,4).Range.Text = "95-52=43"
$t.Cell(13,5).Range.Text = "52+1=53"
$t.Cell(14,1).Range.Text = "40-15=25"
$t.Cell(14,2).Range.Text = "28+8=36"
$t.Cell(14,3).Range.Text = "80-76=4"
$t.Cell(14,4).Range.Text = "1+17=18"
$t.Cell(14,5).Range.Text = "11-7=4"
$t.Cell(15,1).Range.Text = "94-69=25"
$t.Cell(15,2).Range.Text = "3+91=94"
$t.Cell(15,3).Range.Text = "60+22=82"
$t.Cell(15,4).Range.Text = "67-7=60"
$t.Cell(15,5).Range.Text = "33+0=33"
$t.Cell(16,1).Range.Text = "4+90=94"
$t.Cell(16,2).Range.Text = "42+13=55"
$t.Cell(16,3).Range.Text = "26+66=92"
$t.Cell(16,4).Range.Text = "54+45=99"
$t.Cell(16,5).Range.Text = "1+90=91"
$t.Cell(17,1).Range.Text = "13+10=23"
$t.Cell(17,2).Range.Text = "57-56=1"
$t.Cell(17,3).Range.Text = "24+10=34"
$t.Cell(17,4).Range.Text = "15+16=31"
$t.Cell(17,5).Range.Text = "26+5=31"
$t.Cell(18,1).Range.Text = "65-63=2"
$t.Cell(18,2).Range.Text = "36+49=85"
$t.Cell(18,3).Range.Text = "57+22=79"
$t.Cell(18,4).Range.Text = "29+37=66"
$t.Cell(18,5).Range.Text = "10+65=75"
$t.Cell(19,1).Range.Text = "15-11=4"
$t.Cell(19,2).Range.Text = "66+27=93"
$t.Cell(19,3).Range.Text = "41-28=13"
$t.Cell(19,4).Range.Text = "67-37=30"
$t.Cell(19,5).Range.Text = "0+7=7"
$t.Cell(20,1).Range.Text = "73-55=18"
$t.Cell(20,2).Range.Text = "55-37=18"
$t.Cell(20,3).Range.Text = "49+3=52"
$t.Cell(20,4).Range.Text = "50+38=88"
$t.Cell(20,5).Range.Text = "36-15=21"
